# Auto-generated Excel COM-interop script
# Commit: Add data for 2025-03-11
# Updates column L (2025 cumulative) and a few K (2024) correction values
# across the "Citywide Totals", "By Neighborhood", and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 1048
$ws.Range('L3').Value = 1057
$ws.Range('L4').Value = 291
$ws.Range('K5').Value = 586
$ws.Range('L5').Value = 71
$ws.Range('L6').Value = 1081
$ws.Range('K7').Value = 27534
$ws.Range('L7').Value = 3548

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 64
$ws.Range('L3').Value = 67
$ws.Range('L4').Value = 16
$ws.Range('L5').Value = 6
$ws.Range('L6').Value = 64
$ws.Range('L7').Value = 217

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 21
$ws.Range('L3').Value = 38
$ws.Range('L4').Value = 3
$ws.Range('L7').Value = 80

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 32
$ws.Range('L3').Value = 57
$ws.Range('L7').Value = 153

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L2').Value = 18
$ws.Range('L3').Value = 19
$ws.Range('L7').Value = 52

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L3').Value = 31
$ws.Range('K4').Value = 40
$ws.Range('L5').Value = 6
$ws.Range('K7').Value = 910
$ws.Range('L7').Value = 121

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 24
$ws.Range('L4').Value = 13
$ws.Range('L7').Value = 112
$ws.Range('L8').Value = 217
$ws.Range('L11').Value = 57
$ws.Range('L13').Value = 4
$ws.Range('L15').Value = 24
$ws.Range('L16').Value = 9
$ws.Range('L18').Value = 28
$ws.Range('L19').Value = 108
$ws.Range('L20').Value = 90
$ws.Range('L26').Value = 5
$ws.Range('L27').Value = 40
$ws.Range('L29').Value = 166
$ws.Range('L31').Value = 39
$ws.Range('L33').Value = 153
$ws.Range('L36').Value = 61
$ws.Range('K37').Value = 910
$ws.Range('L37').Value = 121
$ws.Range('L42').Value = 111
$ws.Range('L47').Value = 30
$ws.Range('L49').Value = 21
$ws.Range('L52').Value = 66
$ws.Range('L54').Value = 76
$ws.Range('K63').Value = 80
$ws.Range('L63').Value = 12
$ws.Range('L66').Value = 4
$ws.Range('L67').Value = 129
$ws.Range('L76').Value = 42
$ws.Range('L78').Value = 55
$ws.Range('L83').Value = 80
$ws.Range('L84').Value = 37
$ws.Range('L85').Value = 186
$ws.Range('L86').Value = 27
$ws.Range('L88').Value = 53
$ws.Range('L89').Value = 47
$ws.Range('L91').Value = 45
$ws.Range('L94').Value = 43
$ws.Range('L95').Value = 52
$ws.Range('L96').Value = 33
$ws.Range('L97').Value = 42
$ws.Range('L98').Value = 30
$ws.Range('K101').Value = 27534
$ws.Range('L101').Value = 3548

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 39

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 43
$ws.Range('L4').Value = 12
$ws.Range('L7').Value = 129

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L2').Value = 15
$ws.Range('L6').Value = 8
$ws.Range('L7').Value = 37

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('L6').Value = 8
$ws.Range('L7').Value = 21

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L4').Value = 3
$ws.Range('L6').Value = 38
$ws.Range('L7').Value = 76

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 58
$ws.Range('L3').Value = 55
$ws.Range('L4').Value = 7
$ws.Range('L7').Value = 166

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L6').Value = 34
$ws.Range('L7').Value = 108

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L6').Value = 21
$ws.Range('L7').Value = 42

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L3').Value = 24
$ws.Range('L4').Value = 8
$ws.Range('L7').Value = 111

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('L2').Value = 1
$ws.Range('L6').Value = 4

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L2').Value = 18
$ws.Range('L7').Value = 55

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 17
$ws.Range('L3').Value = 5
$ws.Range('L7').Value = 33

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L6').Value = 8
$ws.Range('L7').Value = 45

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 28
$ws.Range('L3').Value = 24
$ws.Range('L4').Value = 5
$ws.Range('L7').Value = 90

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L4').Value = 2
$ws.Range('L7').Value = 28

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L3').Value = 15
$ws.Range('L7').Value = 61

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L3').Value = 42
$ws.Range('L7').Value = 112

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L6').Value = 15
$ws.Range('L7').Value = 43

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L4').Value = 3
$ws.Range('L6').Value = 6
$ws.Range('L7').Value = 30

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L2').Value = 9
$ws.Range('L7').Value = 24

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('L6').Value = 20
$ws.Range('L7').Value = 30

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('L6').Value = 5
$ws.Range('L7').Value = 5

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('L6').Value = 2
$ws.Range('L7').Value = 4

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 20
$ws.Range('L6').Value = 17
$ws.Range('L7').Value = 57

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L3').Value = 8
$ws.Range('L7').Value = 24

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('L2').Value = 5
$ws.Range('L6').Value = 31
$ws.Range('L7').Value = 42

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L3').Value = 16
$ws.Range('L6').Value = 24
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L2').Value = 14
$ws.Range('L7').Value = 47

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L2').Value = 11
$ws.Range('L7').Value = 40

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L3').Value = 4
$ws.Range('L7').Value = 27

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 47
$ws.Range('L3').Value = 82
$ws.Range('L7').Value = 186

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L6').Value = 17
$ws.Range('L7').Value = 66

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('L2').Value = 3
$ws.Range('L7').Value = 13

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('L3').Value = 1
